$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking values in column D so Excel does not
# auto-convert strings like "1.00" or "0.325" into numbers, then restore the
# original (unstyled) cell style so no visual/style diff is introduced.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '67.773.20'
$ws.Range('E2').Value = '  -1.97%  '
Set-TextValue 'D3' '2.390.47'
$ws.Range('E3').Value = '  -3.21%  '
Set-TextValue 'D4' '1.00'
$ws.Range('E4').Value = '  -0.02%  '
Set-TextValue 'D5' '549.13'
$ws.Range('E5').Value = '  -1.93%  '
Set-TextValue 'D6' '156.89'
$ws.Range('E6').Value = '  -3.90%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -0.60%  '
$ws.Range('E9').Value = '  +2.56%  '
$ws.Range('E10').Value = '  -1.58%  '
Set-TextValue 'D11' '0.325'
$ws.Range('E11').Value = '  -3.00%  '
Set-TextValue 'D12' '4.71'
$ws.Range('E12').Value = '  -2.31%  '
Set-TextValue 'D13' '67.693.97'
$ws.Range('E13').Value = '  -1.99%  '
$ws.Range('E14').Value = '  -1.53%  '
Set-TextValue 'D15' '22.68'
$ws.Range('E15').Value = '  -4.00%  '
Set-TextValue 'D16' '10.22'
$ws.Range('E16').Value = '  -5.49%  '
Set-TextValue 'D17' '327.77'
$ws.Range('E17').Value = '  -4.36%  '
Set-TextValue 'D18' '6.71'
$ws.Range('E18').Value = '  -5.44%  '
Set-TextValue 'D19' '3.73'
$ws.Range('E19').Value = '  -2.03%  '
$ws.Range('E20').Value = '  -0.48%  '
$ws.Range('E21').Value = '  -5.20%  '
Set-TextValue 'D22' '65.36'
$ws.Range('E22').Value = '  -2.77%  '
Set-TextValue 'D23' '3.57'
$ws.Range('E23').Value = '  -3.29%  '
$ws.Range('E24').Value = '  -3.43%  '
Set-TextValue 'D25' '0.0₃0786'
$ws.Range('E25').Value = '  -3.89%  '
Set-TextValue 'D26' '6.93'
$ws.Range('E26').Value = '  -3.60%  '
Set-TextValue 'D27' '0.999'
$ws.Range('E27').Value = '  -0.03%  '
Set-TextValue 'D28' '413.41'
$ws.Range('E28').Value = '  -6.27%  '
$ws.Range('E29').Value = '  -3.01%  '
Set-TextValue 'D30' '1.57'
$ws.Range('E30').Value = '  -2.73%  '
Set-TextValue 'D31' '157.24'
$ws.Range('E31').Value = '  +0.71%  '
Set-TextValue 'D32' '18.98'
$ws.Range('E32').Value = '  -0.38%  '
$ws.Range('E33').Value = '  -0.03%  '
Set-TextValue 'D34' '17.55'
$ws.Range('E34').Value = '  -2.08%  '
$ws.Range('E35').Value = '  -4.76%  '
$ws.Range('E36').Value = '  -4.22%  '
$ws.Range('B37').Value = 'Stacks'
$ws.Range('C37').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D37' '1.44'
$ws.Range('E37').Value = '  -2.34%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue 'D38' '4.18'
$ws.Range('E38').Value = '  -6.17%  '
Set-TextValue 'D39' '1.04'
$ws.Range('E39').Value = '  -5.43%  '
Set-TextValue 'D40' '127.40'
$ws.Range('E40').Value = '  -4.46%  '
Set-TextValue 'D41' '3.24'
$ws.Range('E41').Value = '  -3.43%  '
Set-TextValue 'D42' '1.90'
$ws.Range('E42').Value = '  -8.53%  '
Set-TextValue 'D43' '0.0702'
$ws.Range('E43').Value = '  -2.33%  '
$ws.Range('E44').Value = '  -3.05%  '
Set-TextValue 'D45' '0.550'
$ws.Range('E45').Value = '  -2.47%  '
$ws.Range('E46').Value = '  -0.05%  '
$ws.Range('E47').Value = '  -1.17%  '
Set-TextValue 'D48' '1.33'
$ws.Range('E48').Value = '  -8.03%  '
Set-TextValue 'D49' '16.33'
$ws.Range('E49').Value = '  -3.60%  '
Set-TextValue 'D50' '0.0423'
$ws.Range('E50').Value = '  -1.83%  '
Set-TextValue 'D51' '0.0₆0199'
$ws.Range('E51').Value = '  -6.14%  '
